$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'304.30"
$ws.Cells.Item(2, 5).Value = "'0.68%"

$ws.Cells.Item(3, 4).Value = "'35.64"
$ws.Cells.Item(3, 5).Value = "'-4.45%"

$ws.Cells.Item(4, 4).Value = "'5.073"
$ws.Cells.Item(4, 5).Value = "'1.36%"

$ws.Cells.Item(5, 4).Value = "'0.07861"
$ws.Cells.Item(5, 5).Value = "'0.63%"

$ws.Cells.Item(6, 4).Value = "'2.129"
$ws.Cells.Item(6, 5).Value = "'-2.87%"

$ws.Cells.Item(7, 4).Value = "'7.902"
$ws.Cells.Item(7, 5).Value = "'-1.64%"

$ws.Cells.Item(8, 2).Value = "MXToken"
$ws.Cells.Item(8, 3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Cells.Item(8, 4).Value = "'0.9193"
$ws.Cells.Item(8, 5).Value = "'0.72%"

$ws.Cells.Item(9, 2).Value = "LiechtensteinCryptoassetsExchange"
$ws.Cells.Item(9, 3).Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Cells.Item(9, 4).Value = "'0.09752"
$ws.Cells.Item(9, 5).Value = "'0.58%"

$ws.Cells.Item(10, 2).Value = "WazirX"
$ws.Cells.Item(10, 3).Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Cells.Item(10, 4).Value = "'0.1856"
$ws.Cells.Item(10, 5).Value = "'-1.14%"

$ws.Cells.Item(11, 2).Value = "MandalaExchangeToken"
$ws.Cells.Item(11, 3).Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Cells.Item(11, 4).Value = "'0.08638"
$ws.Cells.Item(11, 5).Value = "'-0.51%"

$ws.Cells.Item(12, 2).Value = "BitrueCoin"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Cells.Item(12, 4).Value = "'0.03549"
$ws.Cells.Item(12, 5).Value = "'0.64%"

$ws.Cells.Item(13, 2).Value = "BitMartToken"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Cells.Item(13, 4).Value = "'0.09931"
$ws.Cells.Item(13, 5).Value = "'-0.34%"

$ws.Cells.Item(14, 2).Value = "BitForexToken"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Cells.Item(14, 4).Value = "'0.001429"
$ws.Cells.Item(14, 5).Value = "'-3.24%"

$ws.Cells.Item(15, 2).Value = "TigerCash"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Cells.Item(15, 4).Value = "'0.005632"
$ws.Cells.Item(15, 5).Value = "'-0.32%"

$ws.Cells.Item(16, 2).Value = "LEO"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Cells.Item(16, 4).Value = "'3.459"
$ws.Cells.Item(16, 5).Value = "'0.01%"

$ws.Cells.Item(17, 2).Value = "GateToken"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Cells.Item(17, 4).Value = "'4.102"
$ws.Cells.Item(17, 5).Value = "'1.58%"

$ws.Cells.Item(18, 4).Value = "'2.617"
$ws.Cells.Item(18, 5).Value = "'25.70%"

$ws.Cells.Item(19, 4).Value = "'0.3427"
$ws.Cells.Item(19, 5).Value = "'-1.05%"

$ws.Cells.Item(20, 4).Value = "'5.222"
$ws.Cells.Item(20, 5).Value = "'9.67%"

$ws.Cells.Item(21, 4).Value = "'0.1319"
$ws.Cells.Item(21, 5).Value = "'1.98%"

$ws.Cells.Item(22, 5).Value = "'-0.17%"

$ws.Cells.Item(23, 4).Value = "'0.04551"
$ws.Cells.Item(23, 5).Value = "'-1.96%"

$ws.Cells.Item(24, 4).Value = "'0.005051"
$ws.Cells.Item(24, 5).Value = "'5.35%"

$ws.Cells.Item(25, 4).Value = "'0.001237"
$ws.Cells.Item(25, 5).Value = "'0.53%"

$ws.Cells.Item(27, 4).Value = "'0.0004750"
$ws.Cells.Item(27, 5).Value = "'-0.01%"

$ws.Cells.Item(39, 4).Value = "'0.01848"
$ws.Cells.Item(39, 5).Value = "'5.02%"

$ws.Cells.Item(40, 4).Value = "'0.04723"
$ws.Cells.Item(40, 5).Value = "'-0.30%"

$ws.Cells.Item(41, 4).Value = "'0.007498"
$ws.Cells.Item(41, 5).Value = "'-6.93%"

$ws.Cells.Item(42, 4).Value = "'0.1398"
$ws.Cells.Item(42, 5).Value = "'0.52%"

$ws.Cells.Item(43, 4).Value = "'0.007749"
$ws.Cells.Item(43, 5).Value = "'0.83%"

$ws.Cells.Item(44, 4).Value = "'0.002205"
$ws.Cells.Item(44, 5).Value = "'3.34%"

$ws.Cells.Item(45, 4).Value = "'0.01128"
$ws.Cells.Item(45, 5).Value = "'7.94%"

$ws.Cells.Item(46, 4).Value = "'0.00006337"
$ws.Cells.Item(46, 5).Value = "'4.66%"

$ws.Cells.Item(47, 5).Value = "'-0.01%"

$ws.Cells.Item(48, 4).Value = "'0.0005801"
$ws.Cells.Item(48, 5).Value = "'0.02%"

$ws.Cells.Item(49, 4).Value = "'46.97"
$ws.Cells.Item(49, 5).Value = "'478.07%"

$ws.Cells.Item(50, 4).Value = "'0.002000"
$ws.Cells.Item(50, 5).Value = "'-25.65%"

$ws.Cells.Item(51, 4).Value = "'0.00002101"
$ws.Cells.Item(51, 5).Value = "'-0.01%"
